$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear the cells in row 5 whose style equals their column's default
# style, so they disappear from the saved XML entirely (matches the target).
$ws.Range("A5").Clear()
$ws.Range("C5").Clear()
$ws.Range("E5").Clear()
$ws.Range("I5").Clear()
$ws.Range("J5").Clear()
$ws.Range("K5").Clear()
$ws.Range("L5").Clear()

# The remaining row-5 cells keep their distinct style but lose their value.
$ws.Range("B5:L5").ClearContents()

# Final selection left by the edit.
$ws.Range("E3").Select()
